# Mise a jour du fichier de suivi
# Refresh the "Feuil1" tracking sheet: the list of source filenames in
# column A (rows 2-32) is replaced by the current, alphabetically sorted
# list of XML files (5 placeholder "historiographie_0XX" rows are gone,
# 5 new files were added: the split ed-3-10 volume + the 3 new "rasi"
# files), and every row gets Validation = 1 / Metadonnees = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$files = @(
    @{ Row = 2;  Name = "algarotti_saggio-sopra-l-opera-in-musica_1764.xml" }
    @{ Row = 3;  Name = "arteaga_rivoluzioni_1785.xml" }
    @{ Row = 4;  Name = "calepio_paragone-poesia-tragica_1732.xml" }
    @{ Row = 5;  Name = "martello_della-tragedia-antica-e-moderna_1715.xml" }
    @{ Row = 6;  Name = "napoli-signorelli_addizioni-alla-storia_1798.xml" }
    @{ Row = 7;  Name = "napoli-signorelli_discorso-storico-critico_1783.xml" }
    @{ Row = 8;  Name = "napoli-signorelli_storia-critica-ed-1_1777.xml" }
    @{ Row = 9;  Name = "napoli-signorelli_storia-critica-ed-3-01_1813.xml" }
    @{ Row = 10; Name = "napoli-signorelli_storia-critica-ed-3-02_1813.xml" }
    @{ Row = 11; Name = "napoli-signorelli_storia-critica-ed-3-03_1813.xml" }
    @{ Row = 12; Name = "napoli-signorelli_storia-critica-ed-3-04_1813.xml" }
    @{ Row = 13; Name = "napoli-signorelli_storia-critica-ed-3-05_1813.xml" }
    @{ Row = 14; Name = "napoli-signorelli_storia-critica-ed-3-06_1813.xml" }
    @{ Row = 15; Name = "napoli-signorelli_storia-critica-ed-3-07_1813.xml" }
    @{ Row = 16; Name = "napoli-signorelli_storia-critica-ed-3-08_1813.xml" }
    @{ Row = 17; Name = "napoli-signorelli_storia-critica-ed-3-09_1813.xml" }
    @{ Row = 18; Name = "napoli-signorelli_storia-critica-ed-3-10-1_1813.xml" }
    @{ Row = 19; Name = "napoli-signorelli_storia-critica-ed-3-10-2_1813.xml" }
    @{ Row = 20; Name = "napoli-signorelli_storia-critica-ed2-01_1787.xml" }
    @{ Row = 21; Name = "napoli-signorelli_storia-critica-ed2-02_1787.xml" }
    @{ Row = 22; Name = "napoli-signorelli_storia-critica-ed2-03_1788.xml" }
    @{ Row = 23; Name = "napoli-signorelli_storia-critica-ed2-04_1789.xml" }
    @{ Row = 24; Name = "napoli-signorelli_storia-critica-ed2-05_1789.xml" }
    @{ Row = 25; Name = "napoli-signorelli_storia-critica-ed2-06_1790.xml" }
    @{ Row = 26; Name = "planelli_opera-in-musica_1772.xml" }
    @{ Row = 27; Name = "rasi_comici-italiani-01-01_1897.xml" }
    @{ Row = 28; Name = "rasi_comici-italiani-01-02_1897.xml" }
    @{ Row = 29; Name = "rasi_comici-italiani-02_1897.xml" }
    @{ Row = 30; Name = "riccoboni_observations.xml" }
    @{ Row = 31; Name = "riccoboni_reflexions-historiques-critiques-differents-theatres.xml" }
    @{ Row = 32; Name = "salfi_della-declamazione.xml" }
)

foreach ($file in $files) {
    $ws.Cells.Item($file.Row, 1).Value = $file.Name
    $ws.Cells.Item($file.Row, 2).Value = 1
    $ws.Cells.Item($file.Row, 3).Value = 1
}

# Row 33 keeps its SUM(B2:B32) / SUM(C2:C32) formulas untouched; they
# recalculate automatically to 31 now that every tracked file validates.

# Last selected cell, as left by the author after the edit.
$ws.Range("G31").Select()
